$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $val) {
    $rng = $ws.Range($cellAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "27.625.78"
Set-TextValue "E2" "  -2.46%  "
Set-TextValue "D3" "1.761.39"
Set-TextValue "E3" "  -3.25%  "
Set-TextValue "D4" "1.001"
Set-TextValue "E4" "  -0.08%  "
Set-TextValue "D5" "324.31"
Set-TextValue "E5" "  -1.31%  "
Set-TextValue "D6" "0.9996"
Set-TextValue "E6" "  -0.06%  "
Set-TextValue "D7" "0.4273"
Set-TextValue "E7" "  -1.94%  "
Set-TextValue "D8" "0.3608"
Set-TextValue "E8" "  -2.27%  "
Set-TextValue "D9" "0.07573"
Set-TextValue "E9" "  -1.90%  "
Set-TextValue "D10" "42.32"
Set-TextValue "E10" "  -6.03%  "
Set-TextValue "D11" "1.110"
Set-TextValue "E11" "  -2.79%  "
Set-TextValue "D12" "0.9991"
Set-TextValue "E12" "  -0.14%  "
Set-TextValue "D13" "20.87"
Set-TextValue "E13" "  -6.22%  "
Set-TextValue "D14" "6.079"
Set-TextValue "E14" "  -4.23%  "
Set-TextValue "D15" "7.224"
Set-TextValue "E15" "  -4.75%  "
Set-TextValue "D16" "1.758.29"
Set-TextValue "E16" "  -4.26%  "
Set-TextValue "D17" "93.23"
Set-TextValue "E17" "  -0.40%  "
Set-TextValue "D18" "0.00001071"
Set-TextValue "E18" "  -1.35%  "
Set-TextValue "D19" "0.06403"
Set-TextValue "E19" "  -1.95%  "
Set-TextValue "D20" "0.9992"
Set-TextValue "E20" "  -0.08%  "
Set-TextValue "D21" "17.15"
Set-TextValue "E21" "  -2.34%  "
Set-TextValue "D22" "5.906"
Set-TextValue "E22" "  -6.16%  "
Set-TextValue "D23" "27.673.78"
Set-TextValue "E23" "  -2.43%  "
Set-TextValue "D24" "11.31"
Set-TextValue "E24" "  -3.41%  "
Set-TextValue "D25" "2.124"
Set-TextValue "E25" "  +5.23%  "
Set-TextValue "D26" "162.65"
Set-TextValue "E26" "  +0.74%  "
Set-TextValue "D27" "20.39"
Set-TextValue "E27" "  -2.32%  "
Set-TextValue "D28" "1.959.42"
Set-TextValue "E28" "  -3.89%  "
Set-TextValue "D29" "2.167"
Set-TextValue "E29" "  -6.06%  "
Set-TextValue "D30" "125.67"
Set-TextValue "E30" "  -2.80%  "
Set-TextValue "D31" "1.110"
Set-TextValue "E31" "  -9.00%  "
Set-TextValue "D32" "5.614"
Set-TextValue "E32" "  -6.85%  "
Set-TextValue "D33" "3.659"
Set-TextValue "E33" "  +2.96%  "
Set-TextValue "D34" "0.08918"
Set-TextValue "E34" "  -3.35%  "
Set-TextValue "D35" "12.26"
Set-TextValue "E35" "  -5.83%  "
Set-TextValue "D36" "0.02288"
Set-TextValue "E36" "  -3.24%  "
Set-TextValue "D37" "0.2111"
Set-TextValue "E37" "  -3.49%  "
Set-TextValue "D38" "0.06027"
Set-TextValue "E38" "  -3.07%  "
Set-TextValue "D39" "0.6375"
Set-TextValue "E39" "  -3.71%  "
Set-TextValue "D40" "4.969"
Set-TextValue "E40" "  -4.90%  "
Set-TextValue "D41" "1.190"
Set-TextValue "E41" "  -0.57%  "
Set-TextValue "D42" "0.9987"
Set-TextValue "E42" "  -0.06%  "
Set-TextValue "D43" "1.394"
Set-TextValue "E43" "  -3.27%  "
Set-TextValue "D44" "7.911"
Set-TextValue "E44" "  -3.39%  "
Set-TextValue "D45" "13.42"
Set-TextValue "E45" "  -4.73%  "
Set-TextValue "D46" "0.5954"
Set-TextValue "E46" "  -3.24%  "
Set-TextValue "D47" "3.713"
Set-TextValue "E47" "  -1.24%  "
Set-TextValue "D48" "1.994"
Set-TextValue "E48" "  -2.04%  "
Set-TextValue "D49" "123.50"
Set-TextValue "E49" "  -2.49%  "
Set-TextValue "D50" "1.173"
Set-TextValue "E50" "  +0.94%  "
Set-TextValue "D51" "0.06863"
Set-TextValue "E51" "  -2.33%  "
